# Update the Docker install/verify steps on slide 3 (body placeholder),
# rotating the 4 command lines and tweaking their wording/casing to match
# the published revision of the deck.
#
# NOTE: the COM-interop runtime's TextRange.Text setter preserves runs by
# matching the longest common prefix between the paragraph's existing text
# and the new text, which can (a) leave a stray split run when the prefix
# is non-empty, and (b) misbehave when the new text exactly matches text
# that still lives in a not-yet-updated neighboring paragraph. Writing a
# short, prefix-free placeholder first (so the "old" text for the real
# write shares no prefix with the real new text) avoids both issues and
# yields a single clean run per paragraph, matching the intended XML diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$tr.Paragraphs(4, 1).Text = "#1#"
$tr.Paragraphs(4, 1).Text = "docker --version (Check docker version)"

$tr.Paragraphs(5, 1).Text = "#2#"
$tr.Paragraphs(5, 1).Text = "sudo snap install docker (Other dependencies of other docker packages)"

$tr.Paragraphs(6, 1).Text = "#3#"
$tr.Paragraphs(6, 1).Text = "Sudo systemctl start docker - start Docker "

$tr.Paragraphs(7, 1).Text = "#4#"
$tr.Paragraphs(7, 1).Text = "Sudo systemctl status docker - check the docker running or not"
